$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells holding numeric-looking text must be forced to Text format
# before assignment, otherwise Excel auto-converts them to numbers.
$textCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "27.943.11"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "1.870.92"
$ws.Range("E3").Value = "  -1.80%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").Value = "312.25"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("D7").Value = "0.5022"
$ws.Range("E7").Value = "  -1.39%  "
$ws.Range("D8").Value = "0.3830"
$ws.Range("E8").Value = "  -2.92%  "
$ws.Range("D9").Value = "0.08918"
$ws.Range("E9").Value = "  -7.65%  "
$ws.Range("D10").Value = "1.116"
$ws.Range("E10").Value = "  -1.98%  "
$ws.Range("D11").Value = "41.53"
$ws.Range("E11").Value = "  -1.48%  "
$ws.Range("D12").Value = "6.382"
$ws.Range("E12").Value = "  -0.81%  "
$ws.Range("D13").Value = "20.66"
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").Value = "1.866.52"
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("D15").Value = "7.236"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").Value = "0.9999"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "0.00001098"
$ws.Range("E17").Value = "  -2.50%  "
$ws.Range("D18").Value = "91.07"
$ws.Range("E18").Value = "  -2.12%  "
$ws.Range("D19").Value = "0.06658"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").Value = "18.10"
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").Value = "6.115"
$ws.Range("E22").Value = "  -1.70%  "
$ws.Range("D23").Value = "27.964.02"
$ws.Range("E23").Value = "  -1.39%  "
$ws.Range("D24").Value = "11.49"
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("D25").Value = "2.272"
$ws.Range("E25").Value = "  -1.99%  "
$ws.Range("D26").Value = "2.082.71"
$ws.Range("E26").Value = "  -2.68%  "
$ws.Range("D27").Value = "2.496"
$ws.Range("E27").Value = "  -6.63%  "
$ws.Range("D28").Value = "158.37"
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("D29").Value = "20.65"
$ws.Range("E29").Value = "  -1.82%  "
$ws.Range("D30").Value = "126.22"
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("D31").Value = "0.1060"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").Value = "1.055"
$ws.Range("E32").Value = "  -3.79%  "
$ws.Range("D33").Value = "5.598"
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("D34").Value = "3.599"
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("D35").Value = "9.511"
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("D36").Value = "0.06551"
$ws.Range("E36").Value = "  -1.87%  "
$ws.Range("D37").Value = "0.02400"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("D39").Value = "1.286"
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("D40").Value = "1.206"
$ws.Range("E40").Value = "  -3.35%  "
$ws.Range("D41").Value = "0.6367"
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("D42").Value = "11.49"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("D43").Value = "4.910"
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("D44").Value = "0.9997"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.6006"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "13.16"
$ws.Range("E46").Value = "  -2.81%  "
$ws.Range("D47").Value = "1.278"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").Value = "3.665"
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "1.991"
$ws.Range("E49").Value = "  -2.67%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "1.224"
$ws.Range("E50").Value = "  +2.74%  "
$ws.Range("D51").Value = "120.74"
$ws.Range("E51").Value = "  -2.39%  "
